# Fix first deforestation placeholder data issue.
#
# 1) The four year labels ("2020"/"2021"/"2022"/"2023") on slide 1 are
#    replaced with their generic placeholder tokens ("year1".."year4").
# 2) The first "deforestation" description text box (split across two
#    runs as "d" + "eforestation_text") is corrected/merged into a
#    single run reading "def_description".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Year placeholders -----------------------------------------------
$yearMap = @{
    "Google Shape;57;p13" = "year1"
    "Google Shape;77;p13" = "year2"
    "Google Shape;78;p13" = "year3"
    "Google Shape;79;p13" = "year4"
}

foreach ($shapeName in $yearMap.Keys) {
    $shp = $s.Shapes.Item($shapeName)
    $shp.TextFrame.TextRange.Text = $yearMap[$shapeName]
}

# --- "def_description" textbox ---------------------------------------
# Originally two runs: "d" (lang=en-GB) + "eforestation_text" (lang=en),
# which together render as "deforestation_text". Replace that whole
# (mis-split) word with the corrected single word "def_description" -
# PowerPoint merges it into one run using the first run's properties.
$descShape = $s.Shapes.Item("TextBox 5")
$origWidth = $descShape.Width
$origHeight = $descShape.Height

$descRange = $descShape.TextFrame.TextRange
$descRange.Replace("deforestation_text", "def_description") | Out-Null

# The textbox uses shape auto-fit; keep its footprint exactly as it was
# (the source edit did not resize it).
$descShape.Width = $origWidth
$descShape.Height = $origHeight
